$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 (shifts existing rows 2..60 down to 3..61),
# then clear any inherited formatting so the new row matches the sheet's
# unstyled data rows, and fill in the new day's values.
$ws.Rows(2).Insert()
$ws.Range("A2:D2").ClearFormats()

# Force column A to be stored as text so the date string is not
# reinterpreted as a date serial number.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2026-01-19"
$ws.Range("B2").Value = 783.5
$ws.Range("C2").Value = 1112
$ws.Range("D2").Value = 3610

# Drop the temporary text number format so the new row carries no
# explicit style, same as the rest of the data rows.
$ws.Range("A2:D2").ClearFormats()
